# Reorder the menu rows (A2:D19): move "Chole Bhature" and "Chole chawal"
# (originally rows 2 and 19) down so they sit just before the "Veg Biryani"
# rows, shifting everything else up by two rows. Also move the sheet's
# active selection from D20 to B20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired contents for A2:D19, in row order.
$rows = @(
    @{ A = "Chumin ";                     B = 40;  C = 80;  D = "Chumin.jpeg" },
    @{ A = "Chumin Paneer";                B = 50;  C = 90;  D = "Chumin Paneer.jpeg" },
    @{ A = "Chill Patato";                 B = 40;  C = 80;  D = "chill Patato.jpeg" },
    @{ A = "Honey Chill patato";           B = 50;  C = 90;  D = "Honey Chill patato.jpeg" },
    @{ A = "Aloo tikki Burger ";           B = 0;   C = 40;  D = "Burger Aloo tikki.png" },
    @{ A = "Aloo tikki chees Burger";      B = 0;   C = 50;  D = "XL Burger Aloo tikki chees.jpeg" },
    @{ A = "Paneer tikki chees Burger";    B = 0;   C = 70;  D = "Brioche Bun Burger Paneer tikki chees.jpeg" },
    @{ A = "Tortila Wrap Paneer";          B = 0;   C = 60;  D = "Tortila Wrap Paneer.jpeg" },
    @{ A = "Burrito wrap Paneer";          B = $null; C = 70;  D = "Burrito wrap Paneer.jpeg" },
    @{ A = "Quesadilla Paneer";            B = $null; C = 60;  D = "Quesadilla Paneer.jpeg" },
    @{ A = "Schezwan Grilled Sandwich – Indo-Chinese fusion with Schezwan sauce, veggies, and cheese."; B = $null; C = 80;  D = "Schezwan Grilled Sandwich.png" },
    @{ A = "Club Sandwich (Indian Style) – Multi-layered with veggies, green chutney, Paneer patty, and cheese"; B = $null; C = 100; D = "Multi-layered with veggies Paneer.png" },
    @{ A = "Veg Thail(Chole, Mix Veg, Rice, Raita, salad, 2 Roti, Gulab Zamun)"; B = $null; C = 150; D = "Veg Thali(Chole, Mix Veg, Raita, salad, 4 Roti, Gulab Jamun).png" },
    @{ A = "Veg Special Thail(Sabzi Paneer, Mix Veg, Rice, Raita, salad, 2 Roti, Gulab Zamun)"; B = $null; C = 175; D = "Veg Special Thail(Sabzi Paneer, Mix Veg, Raita, salad, 4 Roti, Gulab Zamun).png" },
    @{ A = "Chole Bhature";                B = 50;  C = 80;  D = "Chole Bhature.jpeg" },
    @{ A = "Chole chawal";                 B = 50;  C = 80;  D = "Chole Chawal.jpeg" },
    @{ A = "Veg Biryani Soya with Garlic Mayo Dip *1";   B = $null; C = 100; D = "Veg Biryani Soya with Garlic Mayo Dip 1.png" },
    @{ A = "Veg Biryani Paneer with Dip Garlic Mayo Dip *1"; B = $null; C = 150; D = "Veg Biryani Paneer with Dip Garlic Mayo Dip 1.png" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    if ($null -eq $row.B) {
        $ws.Cells.Item($r, 2).ClearContents()
    } else {
        $ws.Cells.Item($r, 2).Value = $row.B
    }
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}

# Update the recorded selection to match the saved view (B20 instead of D20).
$ws.Range("B20").Select()
